$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 20063.8
$ws.Range("I38").Value = 20063.8
$ws.Range("K38").Value = 60191.39999999999
$ws.Range("M38").Value = -59819.39999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 367.42105
$ws.Range("I80").Value = 474.72726
$ws.Range("J80").Value = 219.875
$ws.Range("K80").Value = 1424.18178
$ws.Range("L80").Value = 659.625
$ws.Range("M80").Value = -426.1817799999999
$ws.Range("N80").Value = -2655.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 367.42105
$ws.Range("I83").Value = 474.72726
$ws.Range("J83").Value = 219.875
$ws.Range("K83").Value = 4272.54534
$ws.Range("L83").Value = 1978.875
$ws.Range("M83").Value = 719.4546600000003
$ws.Range("N83").Value = -11962.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4173.091
$ws.Range("I116").Value = 3649.75
$ws.Range("J116").Value = 4472.143
$ws.Range("K116").Value = 3649.75
$ws.Range("L116").Value = 4472.143
$ws.Range("M116").Value = -207.75
$ws.Range("N116").Value = -11356.143

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 8390.227999999999
$ws.Range("I132").Value = 9711.117
$ws.Range("K132").Value = 29133.351
$ws.Range("M132").Value = -26603.351

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1476305.6
$ws.Range("I137").Value = 2942180.2
$ws.Range("J137").Value = 10430.883
$ws.Range("K137").Value = 8826540.600000001
$ws.Range("L137").Value = 31292.649
$ws.Range("M137").Value = -8823990.600000001
$ws.Range("N137").Value = -36392.649

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3849.5576
$ws.Range("I32").Value = 3212.442
$ws.Range("J32").Value = 6893.5557
$ws.Range("K32").Value = 3212.442
$ws.Range("L32").Value = 6893.5557
$ws.Range("M32").Value = -2925.442
$ws.Range("N32").Value = -7467.5557

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 25085.723
$ws.Range("I45").Value = 33841.848
$ws.Range("J45").Value = 2319.8
$ws.Range("K45").Value = 33841.848
$ws.Range("L45").Value = 2319.8
$ws.Range("M45").Value = -33464.848
$ws.Range("N45").Value = -3073.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2332.6667
$ws.Range("I61").Value = 1366.6842
$ws.Range("J61").Value = 4626.875
$ws.Range("K61").Value = 1366.6842
$ws.Range("L61").Value = 4626.875
$ws.Range("M61").Value = -1154.6842
$ws.Range("N61").Value = -5050.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3181.5293
$ws.Range("I122").Value = 3292.9375
$ws.Range("K122").Value = 9878.8125
$ws.Range("M122").Value = -7428.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2052.606
$ws.Range("I132").Value = 1349.44
$ws.Range("J132").Value = 4250
$ws.Range("K132").Value = 4048.32
$ws.Range("L132").Value = 12750
$ws.Range("M132").Value = -1518.32
$ws.Range("N132").Value = -17810

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2332.6667
$ws.Range("I136").Value = 1366.6842
$ws.Range("J136").Value = 4626.875
$ws.Range("K136").Value = 4100.0526
$ws.Range("L136").Value = 13880.625
$ws.Range("M136").Value = -1550.0526
$ws.Range("N136").Value = -18980.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H141").Value = 129995
$ws.Range("J141").Value = 129995
$ws.Range("L141").Value = 129995
$ws.Range("N141").Value = -140355

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2631.2727
$ws.Range("I86").Value = 2282.6667
$ws.Range("J86").Value = 4200
$ws.Range("K86").Value = 2282.6667
$ws.Range("L86").Value = 4200
$ws.Range("M86").Value = -1159.6667
$ws.Range("N86").Value = -6446

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2631.2727
$ws.Range("I89").Value = 2282.6667
$ws.Range("J89").Value = 4200
$ws.Range("K89").Value = 11413.3335
$ws.Range("L89").Value = 21000
$ws.Range("M89").Value = -5797.333500000001
$ws.Range("N89").Value = -32232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4419.2046
$ws.Range("I31").Value = 3330.4783
$ws.Range("J31").Value = 5611.619
$ws.Range("K31").Value = 3330.4783
$ws.Range("L31").Value = 5611.619
$ws.Range("M31").Value = -3035.4783
$ws.Range("N31").Value = -6201.619

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4419.2046
$ws.Range("I34").Value = 3330.4783
$ws.Range("J34").Value = 5611.619
$ws.Range("K34").Value = 3330.4783
$ws.Range("L34").Value = 5611.619
$ws.Range("M34").Value = -3128.4783
$ws.Range("N34").Value = -6015.619

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3624.25
$ws.Range("I132").Value = 3749.75
$ws.Range("J132").Value = 3498.75
$ws.Range("K132").Value = 11249.25
$ws.Range("L132").Value = 10496.25
$ws.Range("M132").Value = -8719.25
$ws.Range("N132").Value = -15556.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3189.0588
$ws.Range("I134").Value = 3080.9333
$ws.Range("K134").Value = 9242.7999
$ws.Range("M134").Value = -6707.7999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 291998.16
$ws.Range("I70").Value = 1000000
$ws.Range("J70").Value = 8797.4
$ws.Range("K70").Value = 1000000
$ws.Range("L70").Value = 8797.4
$ws.Range("M70").Value = -999730
$ws.Range("N70").Value = -9337.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 291998.16
$ws.Range("I73").Value = 1000000
$ws.Range("J73").Value = 8797.4
$ws.Range("K73").Value = 1000000
$ws.Range("L73").Value = 8797.4
$ws.Range("M73").Value = -999064
$ws.Range("N73").Value = -10669.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2167.1667
$ws.Range("I102").Value = 2000.7
$ws.Range("J102").Value = 2999.5
$ws.Range("K102").Value = 2000.7
$ws.Range("L102").Value = 2999.5
$ws.Range("M102").Value = -378.7
$ws.Range("N102").Value = -6243.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H112").Value = 100146.5
$ws.Range("J112").Value = 100146.5
$ws.Range("L112").Value = 100146.5
$ws.Range("N112").Value = -102362.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3726.0688
$ws.Range("I122").Value = 2007.8948
$ws.Range("J122").Value = 6990.6
$ws.Range("K122").Value = 6023.6844
$ws.Range("L122").Value = 20971.8
$ws.Range("M122").Value = -3573.6844
$ws.Range("N122").Value = -25871.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3123.2632
$ws.Range("I132").Value = 2595.8572
$ws.Range("J132").Value = 4600
$ws.Range("K132").Value = 7787.571599999999
$ws.Range("L132").Value = 13800
$ws.Range("M132").Value = -5257.571599999999
$ws.Range("N132").Value = -18860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 21230.285
$ws.Range("J56").Value = 21230.285
$ws.Range("L56").Value = 21230.285
$ws.Range("N56").Value = -22658.285

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9499.75
$ws.Range("J62").Value = 9499.75
$ws.Range("L62").Value = 9499.75
$ws.Range("N62").Value = -10747.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 9499.75
$ws.Range("J65").Value = 9499.75
$ws.Range("L65").Value = 47498.75
$ws.Range("N65").Value = -53738.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3143.1738
$ws.Range("I132").Value = 3136.1177
$ws.Range("K132").Value = 9408.3531
$ws.Range("M132").Value = -6878.3531

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 142864960
$ws.Range("I136").Value = 333334080
$ws.Range("K136").Value = 1000002240
$ws.Range("M136").Value = -999999690

Write-Host "Applied 30 row updates across ALC, ARM, BSM, CRP, GSM, WVR sheets"